$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$r = $p1.Range
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Titel"/><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="52"/><w:szCs w:val="52"/></w:rPr><w:t>IMS-T Buchstaben-Jagd</w:t></w:r><w:r><w:rPr><w:sz w:val="52"/><w:szCs w:val="52"/></w:rPr><w:t>!</w:t></w:r><w:r><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Name:___________________</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)
